$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the date series in column A: rows 90..130, each = previous row + 1
for ($r = 90; $r -le 130; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+1"
}

# New data point scraped on day 44177 (row 111)
$ws.Range("B111").Value = 9765
$ws.Range("C111").Formula = "=B111-B79"

# Restore the view state: scrolled down, C112 selected
$ws.Application.ActiveWindow.ScrollRow = 54
$ws.Range("C112").Select()
